$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Capstone 1")

# Mark every "Student Code"-like column (B) and the blank Topic/Description
# placeholder cells with a Text number format *before* writing so that
# numeric-looking strings (student IDs, some with leading zeros) are kept
# as literal text instead of being auto-coerced into numbers by Excel.
$textRange = $ws.Range("B3:B11")
$textRange.NumberFormat = "@"

# --- existing rows: two students swapped places ---
$ws.Cells.Item(3, 2).Value = "3333332123"
$ws.Cells.Item(3, 3).Value = "Viet"
$ws.Cells.Item(3, 4).Value = "Van"

$ws.Cells.Item(5, 2).Value = "3333333221"
$ws.Cells.Item(5, 3).Value = "Quang"
$ws.Cells.Item(5, 4).Value = "Le"

# --- new rows 6..11 ---
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = "2512512112"
$ws.Cells.Item(6, 3).Value = "Huy"
$ws.Cells.Item(6, 4).Value = "Thanh"
$ws.Cells.Item(6, 5).Value = "CMU-TPM"
$ws.Cells.Item(6, 6).Value = "Demo1234"
$ws.Cells.Item(6, 7).Value = " "
$ws.Cells.Item(6, 8).Value = " "

$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = "0921525812"
$ws.Cells.Item(7, 3).Value = "Nguyen"
$ws.Cells.Item(7, 4).Value = "Viet"
$ws.Cells.Item(7, 5).Value = "CMU-TPM"
$ws.Cells.Item(7, 6).Value = "Demo1234"
$ws.Cells.Item(7, 7).Value = " "
$ws.Cells.Item(7, 8).Value = " "

$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = "9316712115"
$ws.Cells.Item(8, 3).Value = "Nguyen"
$ws.Cells.Item(8, 4).Value = "Huy"
$ws.Cells.Item(8, 5).Value = "CMU-TPM"
$ws.Cells.Item(8, 6).Value = "Demo1234"
$ws.Cells.Item(8, 7).Value = " "
$ws.Cells.Item(8, 8).Value = " "

$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = "0921525813"
$ws.Cells.Item(9, 3).Value = "Nguyen"
$ws.Cells.Item(9, 4).Value = "Viet"
$ws.Cells.Item(9, 5).Value = "CMU-TPM"
$ws.Cells.Item(9, 6).Value = "demo2"
$ws.Cells.Item(9, 7).Value = " "
$ws.Cells.Item(9, 8).Value = " "

$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = "2205612617"
$ws.Cells.Item(10, 3).Value = "Viet"
$ws.Cells.Item(10, 4).Value = "Nguyen"
$ws.Cells.Item(10, 5).Value = "CMU-TPM"
$ws.Cells.Item(10, 6).Value = "demo2"
$ws.Cells.Item(10, 7).Value = " "
$ws.Cells.Item(10, 8).Value = " "

$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = "2342151123"
$ws.Cells.Item(11, 3).Value = "Nguyen"
$ws.Cells.Item(11, 4).Value = "Thanh"
$ws.Cells.Item(11, 5).Value = "CMU-TPM"
$ws.Cells.Item(11, 6).Value = "demo2"
$ws.Cells.Item(11, 7).Value = " "
$ws.Cells.Item(11, 8).Value = " "

# Revert the temporary Text format back to General so the written cells
# don't carry an explicit cell style (matches plain shared-string cells
# elsewhere in the sheet).
$textRange.Style = "Normal"
